$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Highlight the "endOfTestData" marker cells in column A (rows 4, 8, 14) with a
# yellow fill - this introduces the new fill/cellXf entries.
$ws.Range("A4").Interior.Color = 65535
$ws.Range("A8").Interior.Color = 65535
$ws.Range("A14").Interior.Color = 65535

# Add a new test-data block (rows 15-18): a header row, two blank rows and a
# closing "endOfTestData" row - mirroring the existing blocks above.

# Row 15: header row, formatted like the other block headers (row 1/5/9).
$null = $ws.Range("A1:D1").Copy()
$null = $ws.Range("A15:D15").PasteSpecial(-4122)
$ws.Range("A15").Value = "validateCreateCustomerAPI"

# Rows 16-17: blank data rows, formatted like the existing blank data cells.
$null = $ws.Range("B4:D4").Copy()
$null = $ws.Range("A16:D16").PasteSpecial(-4122)
$null = $ws.Range("A17:D17").PasteSpecial(-4122)

# Row 18: closing marker row, formatted like the other highlighted
# "endOfTestData" rows.
$null = $ws.Range("A4:D4").Copy()
$null = $ws.Range("A18:D18").PasteSpecial(-4122)
$ws.Range("A18").Value = "endOfTestData"

# Update the active selection to match the edited workbook.
$null = $ws.Range("A13").Select()
